# Small modification to resume:
#  - Trim the leading sentence "I prefer working on small-to-mid-sized
#    teams. " from the SUMMARY paragraph.
#  - Word's automatic "_GoBack" bookmark (last-edit-location) moves from
#    the end of the document to the start of the edited paragraph as a
#    result.

$d = $word.ActiveDocument

# 1. Remove the leading sentence from the summary paragraph.
$d.Content.Find.Execute(
    "I prefer working on small-to-mid-sized teams. I find them",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I find them", 2
) | Out-Null

# 2. Locate the (new) start of that paragraph's remaining text so we can
#    drop the _GoBack bookmark exactly where Word would leave it after
#    editing there.
$target = $d.Content
$target.Find.Execute(
    "I find them to be better",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "", 0
) | Out-Null

# 3. Re-seat the _GoBack bookmark (Word keeps only one instance, moving it
#    to the most recent edit point) at that collapsed position.
$goBackRange = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
